$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Version 3.8.0 -> 3.9.0
$meta.Range("B3").Value = "3.9.0"

# Experimental now has a value "false" -- write it as a real text formula
# result and paste-special the value back so it lands as shared-string text
# (t="s") rather than a native Boolean, matching the source file's output.
$meta.Range("B7").Formula = "=""false"""
$meta.Range("B7").Copy()
$meta.Range("B7").PasteSpecial(-4163)

# Date updated
$meta.Range("B8").Value = "2024-12-02T17:05:26-06:00"

# The three existing "Contact" rows (10, 11, 12) each get a distinct,
# more informative value instead of the generic placeholder text.
$meta.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$meta.Range("B11").Value = "null (iti@ihe.net)"
$meta.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"

# Jurisdiction value updated
$meta.Range("B13").Value = "Global (Whole world)"

$wb.Save()
